$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Find-RowByAccount($account) {
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End([Microsoft.Office.Interop.Excel.XlDirection]::xlUp).Row
    for ($r = 1; $r -le $lastRow; $r++) {
        $v = $ws.Cells.Item($r, 1).Value()
        if ($v -eq $account) {
            return $r
        }
    }
    return -1
}

function Delete-RowByAccount($account) {
    $row = Find-RowByAccount $account
    if ($row -gt 0) {
        $ws.Rows.Item($row).Delete()
    }
}

function Insert-RowBeforeAccount($beforeAccount, $account, $name, $saldo) {
    $row = Find-RowByAccount $beforeAccount
    if ($row -gt 0) {
        $ws.Rows.Item($row).Insert()
        # Force text format on the account-number cell so leading zeros survive
        $ws.Cells.Item($row, 1).NumberFormat = "@"
        $ws.Cells.Item($row, 1).Value = $account
        $ws.Cells.Item($row, 2).Value = $name
        $ws.Cells.Item($row, 3).Value = $saldo
    }
}

# Remove the old 004450724 / ASSAKO row entirely
Delete-RowByAccount "004450724"

# Update the BLUEMETRIX balance
$row = Find-RowByAccount "001761119"
if ($row -gt 0) {
    $ws.Cells.Item($row, 3).Value = 174337.36
}

# Remove the 004751770 / DILSON row entirely
Delete-RowByAccount "004751770"

# Remove the 004222784 / RAFAEL row entirely
Delete-RowByAccount "004222784"

# Remove the 004517080 / TATIANA row entirely
Delete-RowByAccount "004517080"

# Remove the 004267119 / ANA row entirely
Delete-RowByAccount "004267119"

# Remove the 004515341 / BRUNO row entirely
Delete-RowByAccount "004515341"

# Remove the old 005295509 / BHRUNA row (balance 99.54) before re-adding it elsewhere
Delete-RowByAccount "005295509"

# Insert the new 005428871 / ROSANGELA row just before 004693349 / CATARINE
Insert-RowBeforeAccount "004693349" "005428871" "ROSANGELA" 16760.86

# Insert the new (moved) 005295509 / BHRUNA row just before 004392159 / RODRIGO, with updated balance
Insert-RowBeforeAccount "004392159" "005295509" "BHRUNA" 1202.2
